$wb = $excel.ActiveWorkbook

# =====================================================================
# Step 1: Update "总计" (Total) summary sheet - add 2022-Q3 summary row
# =====================================================================
$wsTotal = $wb.Worksheets.Item(1)

# Insert new row 2 (shifts existing data rows 2-8 down to 3-9)
$wsTotal.Range("A2").EntireRow.Insert()

# The Insert operation copies the header row's formatting down into the
# new blank row; strip that back off B2:D2 so they look like normal data
# cells (no special style), matching the other data rows.
$wsTotal.Range("B2:D2").ClearFormats()

# Populate the new row with the 2022-Q3 summary figures
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 39
$wsTotal.Cells.Item(2, 4).Value = 8.41

# Give the new A2 index cell the same style as the other index cells
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

# Renumber the index column for the rows that shifted down (now 1..7)
for ($i = 1; $i -le 7; $i++) {
    $wsTotal.Cells.Item($i + 2, 1).Value = $i
}

# =====================================================================
# Step 2: Insert a brand-new "2022-Q3" worksheet right after "总计",
# containing the per-fund holding breakdown for the new quarter.
# =====================================================================
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

# Header row (row 1), columns B..H
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $wsQ3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Copy the header cell style (bold, centered, bordered) from the "总计"
# sheet's existing B1 header cell, and broadcast it across B1:H1.
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

# Per-fund data rows (rows 2..40)
$data = @"
0	260112	景顺长城能源基建混合A	22.22	61.72	7.60	1.6887	4
1	000979	景顺长城沪港深精选股票	20.32	80.04	8.01	1.6276	3
2	008850	景顺长城价值稳进三年定期开放灵活配置混合	17.29	84.84	7.37	1.2743	6
3	008715	景顺长城价值驱动一年持有期灵活配置混合	8.44	90.91	9.81	0.8280	1
4	009098	景顺长城价值领航两年持有期混合	7.16	90.92	9.80	0.7017	1
5	008060	景顺长城价值边际灵活配置混合A	5.45	80.34	7.96	0.4338	4
6	161611	融通内需驱动混合A/B	8.72	90.68	4.53	0.3950	4
7	502000	西部利得中证500指数增强（LOF）A	17.77	87.95	1.71	0.3039	10
8	014109	融通内需驱动混合C	4.06	90.68	4.53	0.1839	4
9	501059	西部利得中证国有企业红利指数增强（LOF）A	4.10	87.88	2.67	0.1095	4
10	009874	九泰久睿量化股票A	3.15	93.83	3.27	0.1030	9
11	009439	西部利得中证国有企业红利指数增强（LOF）C	3.44	87.88	2.67	0.0918	4
12	009300	西部利得中证500指数增强（LOF）C	4.94	87.95	1.71	0.0845	10
13	002252	融通成长30灵活配置混合A/B	1.56	93.13	4.68	0.0730	5
14	012708	东方红中证东方红红利低波动指数A	3.27	93.80	1.99	0.0651	2
15	015779	景顺长城价值边际灵活配置混合C	0.79	80.34	7.96	0.0629	4
16	014106	融通成长30灵活配置混合C	1.34	93.13	4.68	0.0627	5
17	512890	华泰柏瑞中证红利低波动ETF	2.60	99.50	2.21	0.0575	9
18	008134	鹏华优选价值股票	1.80	92.72	3.02	0.0544	9
19	012879	中信建投量化精选6个月持有期混合C	3.33	90.73	1.23	0.0410	3
20	007499	光大保德信风格轮动混合C	1.63	91.26	1.49	0.0243	10
21	012878	中信建投量化精选6个月持有期混合A	1.67	90.73	1.23	0.0205	3
22	014344	鹏华中证500指数增强A	1.14	92.67	1.66	0.0189	8
23	010120	九泰久福量化股票A	0.54	93.91	3.28	0.0177	10
24	001897	九泰久盛量化先锋灵活配置混合A	0.50	93.59	3.25	0.0162	10
25	009043	九泰久信量化股票	0.43	93.60	3.27	0.0141	9
26	012709	东方红中证东方红红利低波动指数C	0.67	93.80	1.99	0.0133	2
27	014345	鹏华中证500指数增强C	0.73	92.67	1.66	0.0121	8
28	002305	光大保德信风格轮动混合A	0.80	91.26	1.49	0.0119	10
29	004510	九泰久盛量化先锋灵活配置混合C	0.28	93.59	3.25	0.0091	10
30	080015	长盛中小盘精选混合	0.13	84.41	2.94	0.0038	5
31	005260	银华稳健增利灵活配置混合A	0.28	91.67	0.77	0.0022	8
32	006157	财通量化核心优选混合	0.11	83.87	1.64	0.0018	3
33	004546	建信量化优享定期开放灵活配置混合	0.16	25.68	1.01	0.0016	6
34	010121	九泰久福量化股票C	0.04	93.91	3.28	0.0013	10
35	006957	长江量化匠心甄选股票C	0.07	90.77	1.86	0.0013	10
36	006911	长江量化匠心甄选股票A	0.01	90.77	1.86	0.0002	10
37	005261	银华稳健增利灵活配置混合C	0.02	91.67	0.77	0.0002	8
38	016399	九泰久睿量化股票C	0.00	93.83	3.27	0	9
"@

$lines = $data -split "`n"
$r = 2
foreach ($line in $lines) {
    if ($line.Trim().Length -eq 0) { continue }
    $parts = $line -split "`t"

    $idx    = [int]$parts[0]
    $code   = $parts[1]
    $name   = $parts[2]
    $scale  = $parts[3]
    $stockPos = $parts[4]
    $posPct = $parts[5]
    $mktVal = $parts[6]
    $rank   = [int]$parts[7]

    $wsQ3.Cells.Item($r, 1).Value = $idx
    # Fund code / scale / stock-position / position-pct / market-value are
    # stored as *text* in the source data (even though they look numeric),
    # so force them with a leading apostrophe (Excel''s "treat as text"
    # prefix) to avoid the leading zero on the code being stripped.
    $wsQ3.Cells.Item($r, 2).Value = "'" + $code
    $wsQ3.Cells.Item($r, 3).Value = $name
    $wsQ3.Cells.Item($r, 4).Value = "'" + $scale
    $wsQ3.Cells.Item($r, 5).Value = "'" + $stockPos
    $wsQ3.Cells.Item($r, 6).Value = "'" + $posPct
    if ($r -eq 40) {
        # Last row''s market-value column is stored as a plain number (0)
        # in the source data, not as text.
        $wsQ3.Cells.Item($r, 7).Value = 0
    } else {
        $wsQ3.Cells.Item($r, 7).Value = "'" + $mktVal
    }
    $wsQ3.Cells.Item($r, 8).Value = $rank

    $r = $r + 1
}

# Apply the same index-column style used elsewhere (A2:A40)
$wsTotal.Range("A3").Copy()
$wsQ3.Range("A2:A40").PasteSpecial(-4122)

Write-Output "edit complete"
